$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixa")

# --- Update existing rows 2 and 3 IDs ---
$ws.Cells.Item(2, 2).Value = 72481
$ws.Cells.Item(3, 2).Value = 72480

# --- Build new row 4 (copy of row 3's layout/format), values except P/Q/R/S ---
$ws.Cells.Item(4, 1).Value = 67391
$ws.Cells.Item(4, 2).Value = 72479
$ws.Cells.Item(4, 3).Value = $ws.Cells.Item(3, 3).Value()
$ws.Cells.Item(4, 5).Value = $ws.Cells.Item(3, 5).Value()
$ws.Cells.Item(4, 6).Value = $ws.Cells.Item(3, 6).Value()
$ws.Cells.Item(4, 7).Value = $ws.Cells.Item(3, 7).Value()
$ws.Cells.Item(4, 8).Value = $ws.Cells.Item(3, 8).Value()
$ws.Cells.Item(4, 9).Value = $ws.Cells.Item(3, 9).Value()
$ws.Cells.Item(4, 10).Value = $ws.Cells.Item(3, 10).Value()
$ws.Cells.Item(4, 11).Value = $ws.Cells.Item(3, 11).Value()
$ws.Cells.Item(4, 12).Value = $ws.Cells.Item(3, 12).Value()
$ws.Cells.Item(4, 13).Value = $ws.Cells.Item(3, 13).Value()
$ws.Cells.Item(4, 14).Value = $ws.Cells.Item(3, 14).Value()
$ws.Cells.Item(4, 15).Value = $ws.Cells.Item(3, 15).Value()
$ws.Cells.Item(4, 21).Value = $ws.Cells.Item(3, 21).Value()
$ws.Cells.Item(4, 22).Value = $ws.Cells.Item(3, 22).Value()

for ($c = 1; $c -le 22; $c++) {
    $ws.Cells.Item(3, $c).Copy()
    $ws.Cells.Item(4, $c).PasteSpecial(-4122)
}
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(3).RowHeight()

# --- Build new row 5 (copy of row 4's layout/format, column B keeps default style), values except P/Q/R/S ---
$ws.Cells.Item(5, 1).Value = 67392
$ws.Cells.Item(5, 2).Value = 72478
$ws.Cells.Item(5, 3).Value = $ws.Cells.Item(4, 3).Value()
$ws.Cells.Item(5, 5).Value = $ws.Cells.Item(4, 5).Value()
$ws.Cells.Item(5, 6).Value = $ws.Cells.Item(4, 6).Value()
$ws.Cells.Item(5, 7).Value = $ws.Cells.Item(4, 7).Value()
$ws.Cells.Item(5, 8).Value = $ws.Cells.Item(4, 8).Value()
$ws.Cells.Item(5, 9).Value = $ws.Cells.Item(4, 9).Value()
$ws.Cells.Item(5, 10).Value = $ws.Cells.Item(4, 10).Value()
$ws.Cells.Item(5, 11).Value = $ws.Cells.Item(4, 11).Value()
$ws.Cells.Item(5, 12).Value = $ws.Cells.Item(4, 12).Value()
$ws.Cells.Item(5, 13).Value = $ws.Cells.Item(4, 13).Value()
$ws.Cells.Item(5, 14).Value = $ws.Cells.Item(4, 14).Value()
$ws.Cells.Item(5, 15).Value = $ws.Cells.Item(4, 15).Value()
$ws.Cells.Item(5, 21).Value = $ws.Cells.Item(4, 21).Value()
$ws.Cells.Item(5, 22).Value = $ws.Cells.Item(4, 22).Value()

for ($c = 1; $c -le 22; $c++) {
    if ($c -eq 2) { continue }
    $ws.Cells.Item(4, $c).Copy()
    $ws.Cells.Item(5, $c).PasteSpecial(-4122)
}
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight()
$excel.CutCopyMode = 0

# --- Fill in the P/Q/R/S (Hercules/TEC/FOR/DESC) text columns in the exact
#     order the new shared strings were originally authored ---
$ws.Cells.Item(4, 19).Value = "70001DESC"
$ws.Cells.Item(5, 19).Value = "70002DESC"
$ws.Cells.Item(5, 18).Value = "70002FOR"
$ws.Cells.Item(4, 18).Value = "70001FOR"
$ws.Cells.Item(5, 17).Value = "70002TEC"
$ws.Cells.Item(4, 17).Value = "70001TEC"
$ws.Cells.Item(3, 17).Value = "70003TEC"
$ws.Cells.Item(3, 18).Value = "70003FOR"
$ws.Cells.Item(3, 19).Value = "70003DESC"
$ws.Cells.Item(3, 16).Value = "70003Hércules"
$ws.Cells.Item(4, 16).Value = "70001Hércules"
$ws.Cells.Item(5, 16).Value = "70002Hércules"

$ws.Range("R14").Select()
